# Update the cryptos price/volume table (columns D = Price, E = Volume(1h))
# with the latest scraped figures. Cells whose new text parses as a plain
# decimal number (e.g. "500.06") are written with a leading apostrophe so
# Excel keeps them as text instead of silently converting them to numbers -
# matching the original "text-looking-like-a-number" cells in the sheet
# (e.g. "0.999", "499.75"). Values that already contain non-numeric
# characters (multiple dots, percent signs, subscript digits, ...) are
# left as plain literals since Excel cannot parse them as numbers anyway.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.300.32"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "2.275.67"
$ws.Range("E3").Value = "  +1.25%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'500.06"
$ws.Range("E5").Value = "  +1.06%  "
$ws.Range("D6").Value = "'128.85"
$ws.Range("E6").Value = "  +1.27%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("E8").Value = "  -0.48%  "
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("E10").Value = "  +0.54%  "
$ws.Range("E11").Value = "  +3.47%  "
$ws.Range("E12").Value = "  +1.83%  "
$ws.Range("D13").Value = "2.676.84"
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("D14").Value = "'22.67"
$ws.Range("E14").Value = "  +4.43%  "
$ws.Range("D15").Value = "54.257.16"
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("D17").Value = "2.276.90"
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("D18").Value = "'10.23"
$ws.Range("E18").Value = "  +2.18%  "
$ws.Range("E19").Value = "  +1.97%  "
$ws.Range("D20").Value = "'303.42"
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("E21").Value = "  -1.88%  "
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("D23").Value = "'61.21"
$ws.Range("E23").Value = "  -2.85%  "
$ws.Range("D24").Value = "'0.998"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("E25").Value = "  -0.78%  "
$ws.Range("D26").Value = "'7.30"
$ws.Range("E26").Value = "  +2.74%  "
$ws.Range("D27").Value = "'170.40"
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("E28").Value = "  +0.94%  "
$ws.Range("D29").Value = "0.0₃0684"
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("D30").Value = "'5.91"
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("E31").Value = "  +0.64%  "
$ws.Range("E32").Value = "  +0.29%  "
$ws.Range("D33").Value = "'17.75"
$ws.Range("E33").Value = "  +0.77%  "
$ws.Range("D34").Value = "'0.958"
$ws.Range("E34").Value = "  +10.72%  "
$ws.Range("E35").Value = "  +0.33%  "
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("D38").Value = "'0.373"
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("E39").Value = "  +0.37%  "
$ws.Range("E40").Value = "  +0.81%  "
$ws.Range("D41").Value = "'4.82"
$ws.Range("E41").Value = "  -0.69%  "
$ws.Range("D42").Value = "'125.18"
$ws.Range("E42").Value = "  -2.66%  "
$ws.Range("E43").Value = "  +2.03%  "
$ws.Range("D44").Value = "'0.0891"
$ws.Range("E44").Value = "  -0.53%  "
$ws.Range("D45").Value = "'0.545"
$ws.Range("E45").Value = "  -1.02%  "
$ws.Range("D46").Value = "'238.36"
$ws.Range("E46").Value = "  -1.75%  "
$ws.Range("D47").Value = "'0.372"
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("D49").Value = "'10.76"
$ws.Range("E49").Value = "  +0.43%  "
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("E51").Value = "  -0.47%  "
